$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update M3-12mm Bolt quantity from 44 to 42 (row 32)
$ws.Range("C32").Value = 42

# Clear the "M3-10mm nylon Screw" row (row 42) contents in B:C,
# leaving the rest of the row (D/E/F) untouched
$ws.Range("B42:C42").ClearContents()

# Remove the "M3 nylon nut" row entirely (row 43), shifting rows below up
$ws.Rows(43).Delete()

# Update view state to match author's final selection/scroll position
$ws.Range("C40").Select()
$excel.ActiveWindow.ScrollRow = 24
